$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    $cell.NumberFormat = "@"
    $cell.Value = "2012-05-01"
    $cell.Style = "Normal"
}
